$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking values that are stored as plain text
# (coinranking prices, e.g. "260.29" / "0.06100") rather than numbers.
# Flip those cells to text format before writing so Excel keeps the
# exact string (incl. trailing zeros) instead of coercing to a Double,
# then restore the default "Normal" style so no formatting residue is
# left behind.
$priceCells = @("D2","D4","D5","D6","D7","D8","D10","D11","D12","D13","D14","D15","D16","D17","D18","D19","D20","D23","D26","D40","D41","D42","D43","D44","D45","D46","D49")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated cell values in sheet order
$ws.Range("D2").Value = "260.29"
$ws.Range("D4").Value = "6.185"
$ws.Range("D5").Value = "0.06100"
$ws.Range("D6").Value = "6.738"
$ws.Range("D7").Value = "3.485"
$ws.Range("D8").Value = "1.359"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "0.1584"
$ws.Range("E10").Value = "9WazirXWRX"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "0.08052"
$ws.Range("E11").Value = "10MandalaExchangeTokenMDX"
$ws.Range("B12").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C12").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D12").Value = "0.03308"
$ws.Range("E12").Value = "11LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "0.03047"
$ws.Range("E13").Value = "12BitrueCoinBTR"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "0.09304"
$ws.Range("E14").Value = "13BitMartTokenBMX"
$ws.Range("B15").Value = "MCDex"
$ws.Range("C15").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D15").Value = "3.897"
$ws.Range("E15").Value = "14MCDexMCB"
$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").Value = "0.001692"
$ws.Range("E16").Value = "15BitForexTokenBF"
$ws.Range("B17").Value = "CoinExToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D17").Value = "0.04829"
$ws.Range("E17").Value = "16CoinExTokenCET"
$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D18").Value = "0.0006155"
$ws.Range("E18").Value = "17OneONEWorstin24h"
$ws.Range("D19").Value = "0.006191"
$ws.Range("D20").Value = "0.001102"
$ws.Range("D23").Value = "3.693"
$ws.Range("D26").Value = "0.1226"
$ws.Range("D40").Value = "0.04594"
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D41").Value = "0.1119"
$ws.Range("E41").Value = "40BKEXTokenBKK"
$ws.Range("D42").Value = "0.003132"
$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D43").Value = "0.003407"
$ws.Range("E43").Value = "42KickTokenKICK"
$ws.Range("D44").Value = "0.01068"
$ws.Range("D45").Value = "0.002972"
$ws.Range("D46").Value = "0.00005941"
$ws.Range("D49").Value = "0.1139"
$ws.Range("E49").Value = "48BOLOBOLO"

# Restore default styling on the price cells now that the text value is set
foreach ($addr in $priceCells) {
    $ws.Range($addr).Style = "Normal"
}
